# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff has happened for b.md (b.63290e5768f68...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row.
#   B3 (zh-cn status), C3 (de-de status): "Handed back: in sync with en-US" -> "Ready for handoff"
#   D3 (Latest Handoff Date): "2016-27-20 12:27:19" -> "2016-28-20 12:28:15"
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-28-20 12:28:15"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
#   C3 (Status): "Handed back: in sync with en-US" -> "Ready for handoff"
#   D3 (Latest Handoff File): "a.6631f68b...zh-cn.xlf" -> "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
#   E3 (Latest Handoff Datetime): "2016-03-20 12:27:16" -> "2016-03-20 12:28:12"
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-20 12:28:12"

# ---------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row.
#   C3 (Status): "Handed back: in sync with en-US" -> "Ready for handoff"
#   D3 (Latest Handoff File): "a.6631f68b...de-de.xlf" -> "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
#   E3 (Latest Handoff Datetime): "2016-03-20 12:27:19" -> "2016-03-20 12:28:15"
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-20 12:28:15"
